# Refresh the base latency benchmark table (new PDF run results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (aws)
$ws.Range("B2").Value = 41.502507
$ws.Range("C2").Value = 33.943337
$ws.Range("D2").Value = 33.4099445
$ws.Range("E2").Value = 33.283971
$ws.Range("F2").Value = 33.1262525
$ws.Range("G2").Value = 32.935083
$ws.Range("H2").Value = 32.804916
$ws.Range("I2").Value = 32.6085545
$ws.Range("J2").Value = 32.5112955
$ws.Range("K2").Value = 32.409966
$ws.Range("L2").Value = 32.2237045
$ws.Range("M2").Value = 32.163314
$ws.Range("N2").Value = 32.045982
$ws.Range("O2").Value = 32.001113
$ws.Range("P2").Value = 31.9845
$ws.Range("Q2").Value = 31.923474
$ws.Range("R2").Value = 31.9156455
$ws.Range("S2").Value = 31.824707
$ws.Range("T2").Value = 31.848388
$ws.Range("U2").Value = 31.8285285
$ws.Range("V2").Value = 31.830361
$ws.Range("W2").Value = 31.964609
$ws.Range("X2").Value = 32.264929
$ws.Range("Y2").Value = 32.1692005

# Row 3 (cloudflare)
$ws.Range("C3").Value = 18.415742
$ws.Range("D3").Value = 17.990936
$ws.Range("E3").Value = 18.021751
$ws.Range("F3").Value = 18.0168695
$ws.Range("G3").Value = 17.971869
$ws.Range("H3").Value = 17.972097
$ws.Range("I3").Value = 18.106831
$ws.Range("J3").Value = 18.31331
$ws.Range("K3").Value = 18.07888
$ws.Range("L3").Value = 17.866344
$ws.Range("M3").Value = 17.918727
$ws.Range("N3").Value = 18.091247
$ws.Range("O3").Value = 18.502851
$ws.Range("P3").Value = 18.1294
$ws.Range("Q3").Value = 17.817797
$ws.Range("R3").Value = 17.7911055
$ws.Range("S3").Value = 17.690076
$ws.Range("T3").Value = 17.648954
$ws.Range("U3").Value = 17.59945
$ws.Range("V3").Value = 17.596235
$ws.Range("W3").Value = 17.6140005
$ws.Range("X3").Value = 17.6011135
$ws.Range("Y3").Value = 17.107424

# Row 4 (flyio)
$ws.Range("B4").Value = 1049.770767
$ws.Range("C4").Value = 681.5901995
$ws.Range("D4").Value = 41.091653
$ws.Range("E4").Value = 40.576664
$ws.Range("F4").Value = 40.505575
$ws.Range("G4").Value = 40.466164
$ws.Range("H4").Value = 40.442342
$ws.Range("I4").Value = 40.41662
$ws.Range("J4").Value = 40.377234
$ws.Range("K4").Value = 40.368282
$ws.Range("L4").Value = 40.333146
$ws.Range("M4").Value = 40.374788
$ws.Range("N4").Value = 40.3647485
$ws.Range("O4").Value = 40.362897
$ws.Range("P4").Value = 40.34608
$ws.Range("Q4").Value = 40.338186
$ws.Range("R4").Value = 40.338728
$ws.Range("S4").Value = 40.3578675
$ws.Range("T4").Value = 40.329247
$ws.Range("U4").Value = 40.333388
$ws.Range("V4").Value = 40.3309715
$ws.Range("W4").Value = 40.399623
$ws.Range("X4").Value = 40.487854

# Row 5 (google)
$ws.Range("B5").Value = 443.147291
$ws.Range("C5").Value = 65.955372
$ws.Range("D5").Value = 41.767088
$ws.Range("E5").Value = 41.3594285
$ws.Range("F5").Value = 41.159949
$ws.Range("G5").Value = 40.92487850000001
$ws.Range("H5").Value = 40.825615
$ws.Range("I5").Value = 40.7838535
$ws.Range("J5").Value = 40.739068
$ws.Range("K5").Value = 40.713392
$ws.Range("L5").Value = 40.691439
$ws.Range("M5").Value = 40.667528
$ws.Range("N5").Value = 40.6767145
$ws.Range("O5").Value = 40.658665
$ws.Range("P5").Value = 40.7646015
$ws.Range("Q5").Value = 40.994513
$ws.Range("R5").Value = 41.282023
$ws.Range("S5").Value = 41.525326
$ws.Range("T5").Value = 41.632725
$ws.Range("U5").Value = 41.60699
$ws.Range("V5").Value = 41.1323295
$ws.Range("W5").Value = 40.637392
$ws.Range("Y5").Value = 42.758387
